$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(101,9,30,15,60,15),
    @(801,3,67,65,52,45),
    @(1201,2,10,10,10,10),
    @(902,1,0,0,0,0),
    @(901,16,15,45,60,60),
    @(1001,18,30,75,60,72),
    @(401,9,48,67,75,45),
    @(701,3,90,45,97,15),
    @(1203,3,15,15,15,15),
    @(301,6,45,30,60,45),
    @(501,9,52,30,75,45),
    @(601,9,60,67,60,42),
    @(201,9,30,15,45,30),
    @(1202,2,10,10,10,10),
    @(502,0,4,0,0,0),
    @(1101,0,15,30,30,0),
    @(1,0,2,2,2,2),
    @(802,0,4,5,4,0),
    @(2,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(602,0,0,4,0,9),
    @(402,0,0,4,0,0),
    @(702,0,0,0,4,0),
    @(1002,0,0,0,0,9),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
